$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 holds the "id=1" record (col1 of 1 / col2 of 1 / col3 of 1),
# row 6 holds the "id=2" record (col1 of 2 / col2 of 2 / col3 of 2).
# New rows 7..18 repeat those two source rows in this order.
$sourceRowFor = @{
    7  = 4
    8  = 6
    9  = 4
    10 = 4
    11 = 6
    12 = 4
    13 = 4
    14 = 6
    15 = 4
    16 = 4
    17 = 6
    18 = 4
}

foreach ($r in 7..18) {
    $src = $sourceRowFor[$r]
    $ws.Range("B$src`:E$src").Copy()
    $ws.Range("B$r`:E$r").PasteSpecial()
}
